$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format so numeric-looking
# strings (e.g. "319.26") are not auto-converted to floating point numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.763.84'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.479.16'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '319.26'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('E6').Value = '  +2.18%  '
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = '0.0867'
$ws.Range('E10').Value = '  +9.85%  '
$ws.Range('D11').Value = '33.33'
$ws.Range('E11').Value = '  +3.70%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '2.862.05'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '2.488.38'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  +3.05%  '
$ws.Range('D18').Value = '41.741.89'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = '6.48'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').Value = '71.38'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '11.32'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').Value = '240.31'
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('E24').Value = '  +1.80%  '
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '24.80'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '9.84'
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('D30').Value = '36.34'
$ws.Range('E30').Value = '  +3.14%  '
$ws.Range('D31').Value = '158.37'
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').Value = '0.0769'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').Value = '17.43'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').Value = '  +6.86%  '
$ws.Range('E38').Value = '  +3.05%  '
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').Value = '4.02'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('E42').Value = '  +10.67%  '
$ws.Range('D43').Value = '1.987.36'
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('D44').Value = '19.13'
$ws.Range('E44').Value = '  +4.35%  '
$ws.Range('D45').Value = '0.0286'
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').Value = '9.41'
$ws.Range('E47').Value = '  +4.66%  '
$ws.Range('D48').Value = '2.718.56'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('D49').Value = '97.68'
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('E50').Value = '  +4.07%  '
$ws.Range('D51').Value = '67.39'
$ws.Range('E51').Value = '  +0.94%  '

# Restore default cell style so no stray formatting is introduced
# (matches the original workbook, where these cells carry no explicit style).
$dataRange.Style = "Normal"
